$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# Reference cells that already carry the styles we need to reuse:
#  - A4  -> date style (numFmtId 14)
#  - G31 -> time style (numFmtId 20)
$dateStyleSrc = $ws.Range("A4")
$timeStyleSrc = $ws.Range("G31")

function Set-StyledValue($cellRange, $value, $styleSrc) {
    $cellRange.Value = $value
    $styleSrc.Copy()
    $cellRange.PasteSpecial(-4122) | Out-Null
}

# Row 33
Set-StyledValue $ws.Cells.Item(33, 1) 42811 $dateStyleSrc
$ws.Cells.Item(33, 5).Value = "Requirements Specifier"
$ws.Cells.Item(33, 6).Value = "SSD UC6"
Set-StyledValue $ws.Cells.Item(33, 7) 0.34375 $timeStyleSrc
Set-StyledValue $ws.Cells.Item(33, 8) 0.35416666666666669 $timeStyleSrc

# Row 34
$ws.Cells.Item(34, 5).Value = "Any Role"
$ws.Cells.Item(34, 6).Value = "Iterationsplan for iteration 4"
Set-StyledValue $ws.Cells.Item(34, 7) 0.35625000000000001 $timeStyleSrc
Set-StyledValue $ws.Cells.Item(34, 8) 0.39097222222222222 $timeStyleSrc

# Row 35
$ws.Cells.Item(35, 6).Value = "SD og Klassediagram for erVaegtNormal"
Set-StyledValue $ws.Cells.Item(35, 7) 0.40625 $timeStyleSrc
Set-StyledValue $ws.Cells.Item(35, 8) 0.4375 $timeStyleSrc

# Row 36
$ws.Cells.Item(36, 6).Value = "SD og Klassediagram for erArealNormal"
Set-StyledValue $ws.Cells.Item(36, 7) 0.44791666666666669 $timeStyleSrc
Set-StyledValue $ws.Cells.Item(36, 8) 0.47569444444444442 $timeStyleSrc

# Row 37
$ws.Cells.Item(37, 5).Value = "Reviewer"
$ws.Cells.Item(37, 6).Value = "review kode OC5"
Set-StyledValue $ws.Cells.Item(37, 7) 0.51041666666666663 $timeStyleSrc
Set-StyledValue $ws.Cells.Item(37, 8) 0.52083333333333337 $timeStyleSrc

# Row 38
$ws.Cells.Item(38, 5).Value = "Reviewer"
$ws.Cells.Item(38, 6).Value = "review kode OC7"
Set-StyledValue $ws.Cells.Item(38, 7) 0.52083333333333337 $timeStyleSrc
Set-StyledValue $ws.Cells.Item(38, 8) 0.53125 $timeStyleSrc

# Restore the view state: scrolled so row 16 is at the top, with E35 selected
# (matches the edited workbook's sheetView).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E35").Select()
